$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to be treated as literal text so that
# values like "1.000", "0.8850" or "0.00000000120" are not silently
# reinterpreted/normalized as numbers by Excels smart-entry parser.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '29.303.62'
$ws.Range('E2').Value = '  +0.27%  '
Set-TextValue 'D3' '1.874.69'
$ws.Range('E3').Value = '  +0.53%  '
Set-TextValue 'D4' '0.9995'
$ws.Range('E4').Value = '  -0.21%  '
Set-TextValue 'D5' '0.7137'
$ws.Range('E5').Value = '  -0.59%  '
Set-TextValue 'D6' '241.79'
$ws.Range('E6').Value = '  +0.46%  '
Set-TextValue 'D7' '0.9998'
$ws.Range('E7').Value = '  -0.18%  '
Set-TextValue 'D8' '0.3108'
$ws.Range('E8').Value = '  +1.09%  '
Set-TextValue 'D9' '0.07717'
$ws.Range('E9').Value = '  -0.48%  '
Set-TextValue 'D10' '25.08'
$ws.Range('E10').Value = '  +0.19%  '
Set-TextValue 'D11' '0.08382'
$ws.Range('E11').Value = '  +1.55%  '
Set-TextValue 'D12' '1.867.32'
$ws.Range('E12').Value = '  -0.84%  '
Set-TextValue 'D13' '5.216'
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('E14').Value = '  -0.69%  '
Set-TextValue 'D15' '91.34'
$ws.Range('E15').Value = '  +1.15%  '
Set-TextValue 'D16' '29.305.04'
$ws.Range('E16').Value = '  +0.31%  '
Set-TextValue 'D17' '0.000008274'
$ws.Range('E17').Value = '  +6.21%  '
Set-TextValue 'D18' '5.980'
$ws.Range('E18').Value = '  +2.60%  '
Set-TextValue 'D19' '242.57'
$ws.Range('E19').Value = '  -0.24%  '
Set-TextValue 'D20' '2.129.24'
$ws.Range('E20').Value = '  +0.46%  '
$ws.Range('E21').Value = '  +0.59%  '
Set-TextValue 'D22' '0.9991'
$ws.Range('E22').Value = '  -0.28%  '
Set-TextValue 'D23' '7.817'
$ws.Range('E23').Value = '  -1.41%  '
Set-TextValue 'D24' '1.000'
$ws.Range('E24').Value = '  -0.20%  '
Set-TextValue 'D25' '0.1617'
$ws.Range('E25').Value = '  +1.49%  '
Set-TextValue 'D26' '163.19'
$ws.Range('E26').Value = '  +0.67%  '
Set-TextValue 'D27' '9.023'
$ws.Range('E27').Value = '  +1.30%  '
Set-TextValue 'D28' '18.52'
$ws.Range('E28').Value = '  +1.88%  '
Set-TextValue 'D29' '1.505'
$ws.Range('E29').Value = '  +0.79%  '
Set-TextValue 'D30' '4.415'
$ws.Range('E30').Value = '  +1.23%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D31' '1.293'
$ws.Range('E31').Value = '  -1.44%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D32' '4.325'
$ws.Range('E32').Value = '  +5.90%  '
Set-TextValue 'D33' '0.05256'
$ws.Range('E33').Value = '  +1.22%  '
$ws.Range('E34').Value = '  +0.50%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D35' '1.173'
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '0.7471'
$ws.Range('E36').Value = '  +2.63%  '
Set-TextValue 'D37' '2.683'
$ws.Range('E37').Value = '  +0.00%  '
Set-TextValue 'D38' '0.01857'
$ws.Range('E38').Value = '  +0.53%  '
Set-TextValue 'D39' '2.723'
$ws.Range('E39').Value = '  +1.00%  '
Set-TextValue 'D40' '1.156.27'
$ws.Range('E40').Value = '  -0.46%  '
Set-TextValue 'D41' '6.363'
$ws.Range('E41').Value = '  +4.50%  '
Set-TextValue 'D42' '73.03'
$ws.Range('E42').Value = '  +1.26%  '
Set-TextValue 'D43' '0.8850'
$ws.Range('E43').Value = '  -1.91%  '
Set-TextValue 'D44' '105.96'
$ws.Range('E44').Value = '  +4.26%  '
Set-TextValue 'D45' '0.9993'
$ws.Range('E45').Value = '  -0.23%  '
Set-TextValue 'D46' '2.025.28'
$ws.Range('E46').Value = '  +0.39%  '
Set-TextValue 'D47' '1.805'
$ws.Range('E47').Value = '  +2.31%  '
Set-TextValue 'D48' '0.5190'
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '9.378'
$ws.Range('E49').Value = '  +1.29%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D50' '0.00000000120'
$ws.Range('E50').Value = '  +5.12%  '
$ws.Range('E51').Value = '  +1.54%  '
